# Updated cryptos list (Mon Apr 17 15:49:25 UTC 2023 GitHub Actions refresh):
# refresh Price (D) / Volume(1h) (E) figures for every coin row, and fix the
# Algorand/Aptos row ordering (rows 40 & 41 swap Coin/Link/Price/Volume).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells remain text (avoid Excel auto-numeric conversion)
$priceCells = @("D2", "D3", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D46", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.569.21'
$ws.Range("E2").Value = '  -3.54%  '

$ws.Range("D3").Value = '2.087.70'
$ws.Range("E3").Value = '  -1.35%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = '341.77'
$ws.Range("E5").Value = '  -2.07%  '

$ws.Range("E6").Value = '  -0.13%  '

$ws.Range("D7").Value = '0.5130'
$ws.Range("E7").Value = '  -2.56%  '

$ws.Range("D8").Value = '0.4378'
$ws.Range("E8").Value = '  -3.26%  '

$ws.Range("D9").Value = '52.23'
$ws.Range("E9").Value = '  -2.59%  '

$ws.Range("D10").Value = '0.09126'
$ws.Range("E10").Value = '  +1.14%  '

$ws.Range("D11").Value = '1.170'
$ws.Range("E11").Value = '  -0.34%  '

$ws.Range("D12").Value = '24.71'
$ws.Range("E12").Value = '  +0.67%  '

$ws.Range("D13").Value = '2.105.25'
$ws.Range("E13").Value = '  -0.46%  '

$ws.Range("D14").Value = '6.727'
$ws.Range("E14").Value = '  -1.51%  '

$ws.Range("D15").Value = '8.114'
$ws.Range("E15").Value = '  +0.83%  '

$ws.Range("D16").Value = '100.01'
$ws.Range("E16").Value = '  -1.84%  '

$ws.Range("D17").Value = '1.010'
$ws.Range("E17").Value = '  -0.20%  '

$ws.Range("D18").Value = '0.00001143'
$ws.Range("E18").Value = '  -2.20%  '

$ws.Range("D19").Value = '20.88'
$ws.Range("E19").Value = '  +7.64%  '

$ws.Range("D20").Value = '0.06643'
$ws.Range("E20").Value = '  -0.99%  '

$ws.Range("E21").Value = '  -0.16%  '

$ws.Range("D22").Value = '6.142'
$ws.Range("E22").Value = '  -2.63%  '

$ws.Range("D23").Value = '29.602.04'
$ws.Range("E23").Value = '  -3.62%  '

$ws.Range("D24").Value = '12.57'
$ws.Range("E24").Value = '  -2.16%  '

$ws.Range("D25").Value = '2.307'
$ws.Range("E25").Value = '  -3.55%  '

$ws.Range("D26").Value = '2.349.91'
$ws.Range("E26").Value = '  -0.56%  '

$ws.Range("D27").Value = '21.78'
$ws.Range("E27").Value = '  -3.01%  '

$ws.Range("D28").Value = '163.20'
$ws.Range("E28").Value = '  -1.43%  '

$ws.Range("D29").Value = '2.506'
$ws.Range("E29").Value = '  -1.41%  '

$ws.Range("D30").Value = '132.06'
$ws.Range("E30").Value = '  -3.57%  '

$ws.Range("D31").Value = '1.121'
$ws.Range("E31").Value = '  -6.06%  '

$ws.Range("D32").Value = '0.1043'
$ws.Range("E32").Value = '  -3.08%  '

$ws.Range("D33").Value = '1.618'
$ws.Range("E33").Value = '  -1.74%  '

$ws.Range("D34").Value = '6.109'
$ws.Range("E34").Value = '  -4.56%  '

$ws.Range("E35").Value = '  -1.18%  '

$ws.Range("D36").Value = '5.983'
$ws.Range("E36").Value = '  +0.55%  '

$ws.Range("D37").Value = '10.18'
$ws.Range("E37").Value = '  -0.80%  '

$ws.Range("D38").Value = '0.02559'
$ws.Range("E38").Value = '  -3.63%  '

$ws.Range("D39").Value = '0.06663'
$ws.Range("E39").Value = '  -2.81%  '

$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '12.33'
$ws.Range("E40").Value = '  -2.11%  '

$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '0.2219'
$ws.Range("E41").Value = '  -4.31%  '

$ws.Range("D42").Value = '0.6805'
$ws.Range("E42").Value = '  -1.48%  '

$ws.Range("D43").Value = '1.279'
$ws.Range("E43").Value = '  +0.33%  '

$ws.Range("E44").Value = '  +3.15%  '

$ws.Range("E45").Value = '  -4.82%  '

$ws.Range("D46").Value = '2.281'
$ws.Range("E46").Value = '  -2.00%  '

$ws.Range("E47").Value = '  -3.86%  '

$ws.Range("E48").Value = '  -3.07%  '

$ws.Range("D49").Value = '0.00000000335'
$ws.Range("E49").Value = '  -5.43%  '

$ws.Range("D50").Value = '81.27'
$ws.Range("E50").Value = '  -1.41%  '

$ws.Range("D51").Value = '1.171'
$ws.Range("E51").Value = '  -1.87%  '
